$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data for 2022-01-31
$ws.Range("A4").NumberFormat = "d-mmm"
$ws.Range("A4").Value = (Get-Date -Year 2022 -Month 1 -Day 31 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("B4").Value = 1

$ws.Range("C4").Value = "Describing your data"

$ws.Range("C4").Select()
